$wb = $excel.ActiveWorkbook
$wsVariables = $wb.Worksheets.Item(1)
$wsCategories = $wb.Worksheets.Item(2)

# --- sharedStrings text fix: "recetox" cohort renamed to "elspac" ---
$found = $wsCategories.Cells.Find("recetox")
if ($found -ne $null) {
    $found.Value = "elspac"
}

# --- insert a new category row for the new "genrnext" cohort (id 132) ---
$wsCategories.Rows.Item(33).Insert()

$wsCategories.Cells.Item(33, 1).Value = "cohort_id"
$wsCategories.Cells.Item(33, 2).Value = 132
$wsCategories.Cells.Item(33, 3).Value = $false
$wsCategories.Cells.Item(33, 4).Value = "genrnext"

# new row picked up formatting pasted in from elsewhere -> distinct "Normal 2" style
$normal2 = $wb.Styles.Add("Normal 2")
$normal2.Font.Name = "Arial"
$normal2.Font.Size = 10
$wsCategories.Cells.Item(33, 3).Style = "Normal 2"

# --- fix a duplicated urb_area_id code: old row 54 (now row 55) 1803 -> 1804 ---
$wsCategories.Cells.Item(55, 2).Value = 1804

# --- selection / active sheet bookkeeping to mirror the authored edit ---
$wsCategories.Range("A33:XFD33").Select()
$wsCategories.Activate()
